$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 282.7879796666667
$ws.Range("H2").Value = 848.363939
$ws.Range("I2").Value = 0.9674521741401267
$ws.Range("J2").Value = 0.9674521741401266
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.044118333333333
$ws.Range("N2").Value = 6.132354999999999
$ws.Range("O2").Value = 0.1776005292722278
$ws.Range("P2").Value = 0.1776005292722278
$ws.Range("Q2").Value = 578.0520936829272
$ws.Range("R2").Value = 5202.468843146344
$ws.Range("S2").Value = 0.171820018172854
$ws.Range("T2").Value = 0.171820018172854
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 282.7879796666667
$ws.Range("H3").Value = 848.363939
$ws.Range("I3").Value = 0.9674521741401267
$ws.Range("J3").Value = 0.9674521741401266
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.059280333333334
$ws.Range("N3").Value = 21.177841
$ws.Range("O3").Value = 0.6133362746356149
$ws.Range("P3").Value = 0.6133362746356149
$ws.Range("Q3").Value = 1996.279623363967
$ws.Range("R3").Value = 17966.5166102757
$ws.Range("S3").Value = 0.5933735123752315
$ws.Range("T3").Value = 0.5933735123752314
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 282.7879796666667
$ws.Range("H4").Value = 848.363939
$ws.Range("I4").Value = 0.9674521741401267
$ws.Range("J4").Value = 0.9674521741401266
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.406242333333333
$ws.Range("N4").Value = 7.218726999999999
$ws.Range("O4").Value = 0.2090631960921573
$ws.Range("P4").Value = 0.2090631960921573
$ws.Range("Q4").Value = 680.4564080317392
$ws.Range("R4").Value = 6124.107672285652
$ws.Range("S4").Value = 0.2022586435920412
$ws.Range("T4").Value = 0.2022586435920412
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.714696666666668
$ws.Range("H5").Value = 23.14409
$ws.Range("I5").Value = 0.02639291836872237
$ws.Range("J5").Value = 0.02639291836872237
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.044118333333333
$ws.Range("N5").Value = 6.132354999999999
$ws.Range("O5").Value = 0.1776005292722278
$ws.Range("P5").Value = 0.1776005292722278
$ws.Range("Q5").Value = 15.76975289243889
$ws.Range("R5").Value = 141.92777603195
$ws.Range("S5").Value = 0.004687396271323796
$ws.Range("T5").Value = 0.004687396271323796
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 7.714696666666668
$ws.Range("H6").Value = 23.14409
$ws.Range("I6").Value = 0.02639291836872237
$ws.Range("J6").Value = 0.02639291836872237
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.059280333333334
$ws.Range("N6").Value = 21.177841
$ws.Range("O6").Value = 0.6133362746356149
$ws.Range("P6").Value = 0.6133362746356149
$ws.Range("Q6").Value = 54.46020645663223
$ws.Range("R6").Value = 490.1418581096901
$ws.Range("S6").Value = 0.01618773422903407
$ws.Range("T6").Value = 0.01618773422903407
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.714696666666668
$ws.Range("H7").Value = 23.14409
$ws.Range("I7").Value = 0.02639291836872237
$ws.Range("J7").Value = 0.02639291836872237
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.406242333333333
$ws.Range("N7").Value = 7.218726999999999
$ws.Range("O7").Value = 0.2090631960921573
$ws.Range("P7").Value = 0.2090631960921573
$ws.Range("Q7").Value = 18.56342970815889
$ws.Range("R7").Value = 167.07086737343
$ws.Range("S7").Value = 0.005517787868364505
$ws.Range("T7").Value = 0.005517787868364505
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.799090333333333
$ws.Range("H8").Value = 5.397271
$ws.Range("I8").Value = 0.006154907491150983
$ws.Range("J8").Value = 0.006154907491150983
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.044118333333333
$ws.Range("N8").Value = 6.132354999999999
$ws.Range("O8").Value = 0.1776005292722278
$ws.Range("P8").Value = 0.1776005292722278
$ws.Range("Q8").Value = 3.677553533689444
$ws.Range("R8").Value = 33.09798180320499
$ws.Range("S8").Value = 0.001093114828050014
$ws.Range("T8").Value = 0.001093114828050014
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.799090333333333
$ws.Range("H9").Value = 5.397271
$ws.Range("I9").Value = 0.006154907491150983
$ws.Range("J9").Value = 0.006154907491150983
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.059280333333334
$ws.Range("N9").Value = 21.177841
$ws.Range("O9").Value = 0.6133362746356149
$ws.Range("P9").Value = 0.6133362746356149
$ws.Range("Q9").Value = 12.70028300799011
$ws.Range("R9").Value = 114.302547071911
$ws.Range("S9").Value = 0.003775028031349383
$ws.Range("T9").Value = 0.003775028031349383
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.799090333333333
$ws.Range("H10").Value = 5.397271
$ws.Range("I10").Value = 0.006154907491150983
$ws.Range("J10").Value = 0.006154907491150983
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.406242333333333
$ws.Range("N10").Value = 7.218726999999999
$ws.Range("O10").Value = 0.2090631960921573
$ws.Range("P10").Value = 0.2090631960921573
$ws.Range("Q10").Value = 4.329047321557444
$ws.Range("R10").Value = 38.961425894017
$ws.Range("S10").Value = 0.001286764631751586
$ws.Range("T10").Value = 0.001286764631751586
Write-Output "Updated cells for Spp1-Itga9 per Dr Hou advice"
